$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# New key/value pairs to append, matching style used by existing rows
# (column A plain, column B with wrap-text style copied from row 16)
$pairs = @(
    @("victory", "VICTORY"),
    @("combo",   "COMBO"),
    @("bonus",   "BONUS"),
    @("perfect", "PERFECT")
)

$startRow = 17
for ($i = 0; $i -lt $pairs.Length; $i++) {
    $row = $startRow + $i
    $key = $pairs[$i][0]
    $val = $pairs[$i][1]

    $ws.Cells.Item($row, 1).Value = $key
    $ws.Cells.Item($row, 2).Value = $val

    # Match the wrap-text formatting already applied to column B in prior rows
    $ws.Cells.Item($row, 2).WrapText = $true
}

# Update the active selection to reflect the new last row, like Excel would
# after typing down the sheet
$ws.Range("A20").Select()
